$wb = $excel.ActiveWorkbook

# --- Create Tabelle2 as a 2nd sheet (sheetId=3), positioned after Tabelle1 ---
$wsTemp1 = $wb.Worksheets.Add()
$wsTemp2 = $wb.Worksheets.Add()
$wsTemp1fresh = $wb.Worksheets.Item("Sheet1")
$wsTemp1fresh.Delete()
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Name = "Tabelle2"
$ws1fresh = $wb.Worksheets.Item("Tabelle1")
$ws2fresh = $wb.Worksheets.Item("Tabelle2")
$ws2fresh.Move($null, $ws1fresh)

$ws1 = $wb.Worksheets.Item("Tabelle1")
$ws2 = $wb.Worksheets.Item("Tabelle2")

# --- Populate formats by copying from the matching styled cells in Tabelle1 ---
$ws1.Range("D7").Copy()
$ws2.Range("D9","D13","B16","C18","C19","B21").PasteSpecial(-4122)
$ws1.Range("E7").Copy()
$ws2.Range("E9","E13","C16","D18","D19","E20","C21").PasteSpecial(-4122)
$ws1.Range("B16").Copy()
$ws2.Range("B22","C22","D22","E22").PasteSpecial(-4122)
$ws1.Range("A1").Copy()
$ws2.Range("A1").PasteSpecial(-4122)
$ws1.Range("E3").Copy()
$ws2.Range("E3").PasteSpecial(-4122)
$ws1.Range("C3").Copy()
$ws2.Range("C3").PasteSpecial(-4122)
$ws1.Range("B3").Copy()
$ws2.Range("B3","D3","F3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Border for D20 (left+bottom thin) - new style, not present in Tabelle1 ---
$ws2.Range("D20").Borders.Item(7).LineStyle = 1
$ws2.Range("D20").Borders.Item(9).LineStyle = 1

# --- Values / formulas ---
$ws2.Range("A1").Value = "Multiplikation 2er unsigned Int64:"
$ws2.Range("B2").Value = "2 * 4-stellige Zahlen"
$ws2.Range("B3").Value = 5584
$ws2.Range("C3").Value = "*"
$ws2.Range("D3").Value = 4927
$ws2.Range("E3").Value = "="
$ws2.Range("F3").Formula = "=B3*D3"
$ws2.Range("A6").Value = "Jede Zelle entspricht 32-Bit:"
$ws2.Range("A7").Value = "Jeder Kasten entspricht 64-Bit:"
$ws2.Range("D8").Value = "[esp+8]"
$ws2.Range("E8").Value = "[esp+4]"
$ws2.Range("D9").Formula = "=(B3-E9)/100"
$ws2.Range("E9").Formula = "=MOD(B3, 100)"
$ws2.Range("F11").Formula = "=(D9*100+E9)*(D13*100+E13)"
$ws2.Range("G11").Formula = "=IF(F11=F3, `"OK`", `"Achtung`")"
$ws2.Range("D12").Value = "[esp+16]"
$ws2.Range("E12").Value = "[esp+12]"
$ws2.Range("D13").Formula = "=(D3-E13)/100"
$ws2.Range("E13").Formula = "=MOD(D3, 100)"
$ws2.Range("B16").Formula = "=(D13*D9-C16)/100"
$ws2.Range("C16").Formula = "=MOD(D13*D9,100)"
$ws2.Range("C18").Formula = "=(D9*E13-D18)/100"
$ws2.Range("D18").Formula = "=MOD(D9*E13,100)"
$ws2.Range("C19").Formula = "=(D13*E9-D19)/100"
$ws2.Range("D19").Formula = "=MOD(D13*E9,100)"
$ws2.Range("D20").Formula = "=(E9*E13 - E20)/100"
$ws2.Range("E20").Formula = "=MOD(E9*E13, 100)"
$ws2.Range("B21").Formula = "=(SUM(C16:C21)-C23)/100"
$ws2.Range("C21").Formula = "=(SUM(D16:D21)-D23)/100"
$ws2.Range("B23").Formula = "=SUM(B16:B21)"
$ws2.Range("C23").Formula = "=MOD(SUM(C16:C21), 100)"
$ws2.Range("D23").Formula = "=MOD(SUM(D16:D21), 100)"
$ws2.Range("E23").Formula = "=E20"
$ws2.Range("F23").Formula = "=(B23*1000000) + (C23* 10000) + (D23*100) + E23"
$ws2.Range("G23").Formula = "=IF(F23=F11, `"OK`", `"Achtung`")"

# --- Sheet view / selection ---
$ws1.Range("B15").Select()
$ws2.Activate()
$ws2.Range("B21:C21").Select()
